# Add a new game record (row 22) to the Game_Record sheet.
# This mirrors the author appending one more played game on 2026-02-05:
#   Date 46058, 1st=SimpleJack, 2nd=DrSystomatix, 3rd=Doanage, 4th=SiderFace,
#   Start 0.87777777777777777, End 0.95763888888888893
# Everything downstream (H2 COUNTA, A-column ROW()-1 fill, and all the
# Stat_Sheet COUNTIF/ratio formulas) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Game_Record")

$ws.Range("A22").Formula = "=ROW()-1"
$ws.Range("B22").Value = 46058
$ws.Range("C22").Value = "SimpleJack"
$ws.Range("D22").Value = "DrSystomatix"
$ws.Range("E22").Value = "Doanage"
$ws.Range("F22").Value = "SiderFace"
$ws.Range("I22").Value = 0.87777777777777777
$ws.Range("J22").Value = 0.95763888888888893

# Match the Start/End time-of-day formatting used by the row above.
$ws.Range("I22").NumberFormat = $ws.Range("I21").NumberFormat
$ws.Range("J22").NumberFormat = $ws.Range("J21").NumberFormat

# The author's next click landed on the first empty row beneath the new entry.
[void]$ws.Range("A23").Select()
